# Refined metadata to be additional tab
#
# 1. Update the "time_taken" (column F) timestamps on the "data" sheet.
# 2. Add a new "metadata" worksheet (right after "data") that summarizes
#    the panel query that produced the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh the time_taken timestamps on the data sheet ---------------
$dataSheet.Range("F2").Value = "2021-10-05 14:33:16.300511"
$dataSheet.Range("F3").Value = "2021-10-05 14:33:16.300526"
$dataSheet.Range("F4").Value = "2021-10-05 14:33:16.300531"
$dataSheet.Range("F5").Value = "2021-10-05 14:33:16.300536"
$dataSheet.Range("F6").Value = "2021-10-05 14:33:16.300540"
$dataSheet.Range("F7").Value = "2021-10-05 14:33:16.300543"
$dataSheet.Range("F8").Value = "2021-10-05 14:33:16.300546"
$dataSheet.Range("F9").Value = "2021-10-05 14:33:16.300548"
$dataSheet.Range("F10").Value = "2021-10-05 14:33:16.300552"

# --- 2. Add the "metadata" worksheet, placed right after "data" -----------
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Reuse the same header / index cell styling used on the "data" sheet
# (bold font, thin border, centered) instead of re-creating it, so the
# workbook's style table stays as close as possible to the original.
$dataSheet.Range("B1:F1").Copy()
$metaSheet.Range("B1:F1").PasteSpecial(-4122)
$dataSheet.Range("F1").Copy()
$metaSheet.Range("G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$metaSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"

# Data row
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("B2").Value = "Bartter Syndrome"
$metaSheet.Range("C2").Value = 52

# Force "data_version" to stay a text value ("0.17") rather than being
# auto-coerced to the number 0.17.
$metaSheet.Range("D2").NumberFormat = "@"
$metaSheet.Range("D2").Value = "0.17"

$metaSheet.Range("E2").Value = "2021-02-09T09:52:35.102175Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:33:16.294497"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/52/?format=json"

# Keep "data" as the active sheet/tab, same as before the edit.
$dataSheet.Activate() | Out-Null
